# "semana 32 de 2025" - add week-32 column (AI) to the IRA/UCI weekly revision sheet.
# Mirrors the existing week columns (D..AH hold weeks 1..31): the header row gets a
# new label "32" in AI1, and every data row that already reports a week-31 figure in
# AH gets a matching week-32 figure in AI. Rows with no AH value (facilities that
# aren't tracked that far) are left untouched, exactly like the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AI = column 35 (A=1 ... Z=26, AA=27 ... AH=34, AI=35)
$col = 35

# Header label for the new week column.
$ws.Cells.Item(1, $col).Value = 32

# Row -> week-32 value, taken from the commit's new AI column.
$weekValues = [ordered]@{
    2  = 0;
    3  = 0;
    5  = 0;
    6  = 4;
    7  = 0;
    8  = 0;
    9  = 0;
    10 = 0;
    13 = 0;
    15 = 0;
    16 = 0;
    17 = 0;
    23 = 0;
    25 = 0;
    27 = 0;
    28 = 10;
    29 = 1;
    30 = 3;
    31 = 0;
    32 = 0;
    34 = 0;
    35 = 0;
    36 = 0;
    37 = 0;
    38 = 0;
    40 = 0;
    41 = 0;
    42 = 0;
    43 = 0;
    45 = 0;
    46 = 0;
    47 = 0;
    48 = 0;
    49 = 0;
    50 = 0;
    51 = 0;
    53 = 0;
    54 = 0;
    55 = 0;
    56 = 0;
    57 = 0;
    58 = 0
}

foreach ($row in $weekValues.Keys) {
    $ws.Cells.Item($row, $col).Value = $weekValues[$row]
}
